# Append a new login record (test6@gmail.com / test6) as row 7,
# matching the existing username/password columns in the "login" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "test6@gmail.com"
$ws.Range("B7").Value = "test6"
